# Add a new localization-status row for a1c53372-a228-4bad-b83b-7164ea0a7679.md
# across the Overview, zh-cn and de-de sheets ("Generate Report for Handoff").

$wb = $excel.ActiveWorkbook

$fileName      = "a1c53372-a228-4bad-b83b-7164ea0a7679.md"
$pathAndName   = "e2e\a1c53372-a228-4bad-b83b-7164ea0a7679.md"
$status        = "Ready for handoff"
$zhTargetFile  = "a1c53372-a228-4bad-b83b-7164ea0a7679.99a121122415cfe418d5024b564193379b233f88.zh-cn.xlf"
$zhHandoffDate = "2016-09-04 06:47:54"
$deTargetFile  = "a1c53372-a228-4bad-b83b-7164ea0a7679.99a121122415cfe418d5024b564193379b233f88.de-de.xlf"
$deHandoffDate = "2016-09-04 06:47:59"
$overviewDate  = "2016-09-04 06:47:59"
$hrefBase      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/636e15234b0b47723ef6d53d013e59e5485372ff/e2e/a1c53372-a228-4bad-b83b-7164ea0a7679.md"

# ---------------- Overview sheet ----------------
$wsOverview = $wb.Worksheets.Item("Overview")
$rowOverview = 9

$wsOverview.Cells.Item($rowOverview, 1).Value = $fileName
$wsOverview.Cells.Item($rowOverview, 2).Value = $pathAndName
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($rowOverview, 2), $hrefBase, "", "", $pathAndName) | Out-Null
$wsOverview.Cells.Item($rowOverview, 3).Value = ".md"
$wsOverview.Cells.Item($rowOverview, 4).Value = ""
$wsOverview.Cells.Item($rowOverview, 5).Value = $status
$wsOverview.Cells.Item($rowOverview, 6).Value = $status
$wsOverview.Cells.Item($rowOverview, 7).Value = $overviewDate

# ---------------- zh-cn sheet ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$rowZh = 9

$wsZh.Cells.Item($rowZh, 1).Value = $fileName
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rowZh, 1), $hrefBase, "", "", $fileName) | Out-Null
$wsZh.Cells.Item($rowZh, 2).Value = ".md"
$wsZh.Cells.Item($rowZh, 3).Value = $status
$wsZh.Cells.Item($rowZh, 4).Value = "e2e"
$wsZh.Cells.Item($rowZh, 5).Value = "ht"
$wsZh.Cells.Item($rowZh, 6).Value = "False"
$wsZh.Cells.Item($rowZh, 7).Value = $zhTargetFile
$wsZh.Cells.Item($rowZh, 8).Value = $zhHandoffDate
$wsZh.Cells.Item($rowZh, 9).Value = ""
$wsZh.Cells.Item($rowZh, 10).Value = ""
$wsZh.Cells.Item($rowZh, 11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item($rowZh, 12).Value = ""
$wsZh.Cells.Item($rowZh, 13).Value = "True"
$wsZh.Cells.Item($rowZh, 14).Value = ""
$wsZh.Cells.Item($rowZh, 15).Value = "False"
$wsZh.Cells.Item($rowZh, 16).Value = ""

# ---------------- de-de sheet ----------------
$wsDe = $wb.Worksheets.Item("de-de")
$rowDe = 9

$wsDe.Cells.Item($rowDe, 1).Value = $fileName
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rowDe, 1), $hrefBase, "", "", $fileName) | Out-Null
$wsDe.Cells.Item($rowDe, 2).Value = ".md"
$wsDe.Cells.Item($rowDe, 3).Value = $status
$wsDe.Cells.Item($rowDe, 4).Value = "e2e"
$wsDe.Cells.Item($rowDe, 5).Value = "ht"
$wsDe.Cells.Item($rowDe, 6).Value = "False"
$wsDe.Cells.Item($rowDe, 7).Value = $deTargetFile
$wsDe.Cells.Item($rowDe, 8).Value = $overviewDate
$wsDe.Cells.Item($rowDe, 9).Value = ""
$wsDe.Cells.Item($rowDe, 10).Value = ""
$wsDe.Cells.Item($rowDe, 11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item($rowDe, 12).Value = ""
$wsDe.Cells.Item($rowDe, 13).Value = "True"
$wsDe.Cells.Item($rowDe, 14).Value = ""
$wsDe.Cells.Item($rowDe, 15).Value = "False"
$wsDe.Cells.Item($rowDe, 16).Value = ""

# ---------------- Table ranges ----------------
$wb.Worksheets.Item("Overview").ListObjects.Item(1).Resize($wsOverview.Range("A1:G9"))
$wb.Worksheets.Item("zh-cn").ListObjects.Item(1).Resize($wsZh.Range("A1:P9"))
$wb.Worksheets.Item("de-de").ListObjects.Item(1).Resize($wsDe.Range("A1:P9"))
